$d = $word.ActiveDocument

# Helper: replace a whole paragraph's run content with a fresh list of runs
# (each becoming its own <w:r><w:t xml:space="preserve">...</w:t></w:r>).
# We splice OOXML into the paragraph's text range (excluding the trailing
# paragraph mark) via InsertXML so the <w:pPr>/pStyle of the paragraph is
# left untouched while we get full control of run boundaries - something a
# plain Find/Replace cannot do when a paragraph needs to end up with more
# than one run.
function Set-ParagraphRuns($paragraph, [string[]]$texts) {
    $full = $paragraph.Range
    $textRange = $d.Range($full.Start, $full.End - 1)

    $runsXml = ($texts | ForEach-Object {
        $escaped = $_ -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
        "<w:r><w:t xml:space=`"preserve`">$escaped</w:t></w:r>"
    }) -join ''

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $textRange.InsertXML($xml)
}

# Helper: find a paragraph by its exact (trailing-mark-trimmed) text.
function Find-ParagraphByText([string]$text) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.TrimEnd() -eq $text) {
            return $p
        }
    }
    return $null
}

# 1) "Another test just in case" -> "This one comes from the deploy keys"
$p1 = Find-ParagraphByText "Another test just in case"
Set-ParagraphRuns $p1 @("This one comes from the deploy keys")

# 2) "Another test with the right branch." -> "And another one from the deploy keys repo"
$p2 = Find-ParagraphByText "Another test with the right branch."
Set-ParagraphRuns $p2 @("And another one from the deploy keys repo")

# 3) "Fix / else / then is tested" -> three runs:
#    "SSH_DEPLOY_KEY updated (now without a new line at the end)" + " " +
#    "Now using the correct deploy key"
$p3 = Find-ParagraphByText "Fix / else / then is tested"
Set-ParagraphRuns $p3 @(
    "SSH_DEPLOY_KEY updated (now without a new line at the end)",
    " ",
    "Now using the correct deploy key"
)
